$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the H1 title
#    ("Play 123 Boom! Slot Game for Free - Review"). The new paragraph
#    has a bold "Meta description" run followed by a plain run with the
#    description text.
# ---------------------------------------------------------------------
$titlePara = $d.Paragraphs(1)
$titlePara.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs(2)
$metaPara.Style = "Normal"

$metaRange = $metaPara.Range
$metaLabel = "Meta description"
$metaRange.Text = $metaLabel

$boldRange = $d.Range($metaPara.Range.Start, $metaPara.Range.Start + $metaLabel.Length)
$boldRange.Font.Bold = 1

$metaPara.Range.InsertAfter(": Our review of 123 Boom! slot game covers its unique gameplay, symbol design, and bonus features. Try it for free and enjoy explosive wins.")

$plainRange = $d.Range($metaPara.Range.Start + $metaLabel.Length, $metaPara.Range.End)
$plainRange.Font.Bold = 0

# ---------------------------------------------------------------------
# 2) Near the end of the document there used to be a duplicate of the
#    title (bold run) immediately followed by the meta-description
#    paragraph (italic run). The title duplicate is removed entirely,
#    and the italic paragraph's text is replaced with a new AI image
#    generation prompt (formatting/run stays untouched).
# ---------------------------------------------------------------------
$searchRange = $d.Range($titlePara.Range.End, $d.Content.End)
$found = $searchRange.Find.Execute("Play 123 Boom! Slot Game for Free - Review", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$dupTitlePara = $searchRange.Paragraphs(1)
$dupTitlePara.Range.Delete()

$imagePromptPara = $d.Paragraphs($d.Paragraphs.Count)
$imageRange = $d.Range($imagePromptPara.Range.Start, $imagePromptPara.Range.End - 1)
$imageRange.Text = "Create a feature image for ""123 Boom!"" that features a happy Maya warrior with glasses in a cartoon style. In the center of the image, the Maya warrior should be holding a treasure chest overflowing with coins. The Maya warrior should be surrounded by pirate-themed symbols such as crossed swords, black flags, guns, and maps. In the background, there should be a naval battle scene with two ships blasting each other with cannons. The overall style of the image should be fun and adventurous, capturing the spirit of the game."

Write-Output "done"
